# Generate Report for Handback
#
# Row 6 in both the "zh-cn" and "de-de" sheets corresponds to the
# 3884ab68-39a9-49b7-9244-2758b50b27be.md handback file. A handback report
# was generated for it: the target file / handback file / handback datetime
# columns get populated, and because the handback turned out to be stale
# (not built from the latest handoff) an error message is recorded too.

$wb = $excel.ActiveWorkbook

$errorMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0b7fb5071a11e73ea1d45e751355a38a4ce853e6/e2e/3884ab68-39a9-49b7-9244-2758b50b27be.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/29f20441e137c708004d965967dc8480aa47e2d5/e2e/3884ab68-39a9-49b7-9244-2758b50b27be.md."
$latestMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/29f20441e137c708004d965967dc8480aa47e2d5/e2e/3884ab68-39a9-49b7-9244-2758b50b27be.md"
$mdDisplay = "3884ab68-39a9-49b7-9244-2758b50b27be.md"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Widen the "Error Detail" column (P, the 16th column) so the message fits,
# matching the width already used on the wide columns (A, G, I, J).
$wsZh.Columns.Item(16).ColumnWidth = $wsZh.Columns.Item(1).ColumnWidth

# I6 "Latest Target File": add the hyperlink to the handed-back md file.
$wsZh.Hyperlinks.Add($wsZh.Range("I6"), $latestMdUrl, "", "", $mdDisplay)
$wsZh.Range("I6").Font.Underline = $true
$wsZh.Range("I6").Font.Color = 15570276

# J6 "Latest Handback File": same handoff xlf that was handed back.
$wsZh.Range("J6").Value = $wsZh.Range("G6").Value2

# K6 "Latest Handback DateTime".
$wsZh.Range("K6").Value = "2016-08-19 02:40:45"

# P6 "Error Detail".
$wsZh.Range("P6").Value = $errorMessage

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Columns.Item(16).ColumnWidth = $wsDe.Columns.Item(1).ColumnWidth

$wsDe.Hyperlinks.Add($wsDe.Range("I6"), $latestMdUrl, "", "", $mdDisplay)
$wsDe.Range("I6").Font.Underline = $true
$wsDe.Range("I6").Font.Color = 15570276

$wsDe.Range("J6").Value = $wsDe.Range("G6").Value2

$wsDe.Range("K6").Value = "2016-08-19 02:40:52"

$wsDe.Range("P6").Value = $errorMessage
